$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.955.53'
$ws.Range('E2').Value = '  -6.85%  '
$ws.Range('D3').Value = '2.548.82'
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''298.07'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').Value = '''92.13'
$ws.Range('E6').Value = '  -7.03%  '
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -5.83%  '
$ws.Range('D10').Value = '''35.77'
$ws.Range('E10').Value = '  -8.47%  '
$ws.Range('E11').Value = '  -4.29%  '
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('D13').Value = '''0.108'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '2.937.84'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').Value = '2.521.58'
$ws.Range('E15').Value = '  -2.80%  '
$ws.Range('E16').Value = '  -5.48%  '
$ws.Range('D17').Value = '''14.15'
$ws.Range('E17').Value = '  -4.78%  '
$ws.Range('D18').Value = '42.943.51'
$ws.Range('E18').Value = '  -7.07%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''6.65'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0976'
$ws.Range('E20').Value = '  -4.01%  '
$ws.Range('D21').Value = '''12.57'
$ws.Range('E21').Value = '  -1.81%  '
$ws.Range('D22').Value = '''72.13'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').Value = '''260.62'
$ws.Range('E23').Value = '  -11.30%  '
$ws.Range('D24').Value = '''2.91'
$ws.Range('E24').Value = '  -5.08%  '
$ws.Range('D25').Value = '''29.48'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '''2.13'
$ws.Range('E26').Value = '  -4.47%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('E28').Value = '  -7.26%  '
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('D30').Value = '''36.46'
$ws.Range('E30').Value = '  -6.21%  '
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('D32').Value = '''152.15'
$ws.Range('E32').Value = '  -2.43%  '
$ws.Range('D33').Value = '''2.17'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  -5.48%  '
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('E37').Value = '  -6.71%  '
$ws.Range('D38').Value = '''24.14'
$ws.Range('E38').Value = '  +14.57%  '
$ws.Range('E39').Value = '  -3.43%  '
$ws.Range('D40').Value = '''16.62'
$ws.Range('E40').Value = '  +5.58%  '
$ws.Range('D41').Value = '''3.45'
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('E42').Value = '  -6.24%  '
$ws.Range('D43').Value = '''3.83'
$ws.Range('E43').Value = '  -4.32%  '
$ws.Range('D44').Value = '2.077.38'
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('D45').Value = '''0.999'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''85.16'
$ws.Range('E46').Value = '  -13.11%  '
$ws.Range('E47').Value = '  +3.15%  '
$ws.Range('D48').Value = '2.793.24'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''104.30'
$ws.Range('E49').Value = '  -4.06%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''1.70'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').Value = '''8.66'
$ws.Range('E51').Value = '  -8.67%  '
